$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header row's formatting down to the new data row (reuses the
# same cell style index instead of minting a new one).
$ws.Range("A1:K1").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)

$ws.Range("A2").Value = "MentalHealthTreatmentFollowupStatus"
$ws.Range("B2").Value = "Mental Health Treatment Follow-up Status"
$ws.Range("C2").Value = "null#social-history"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G2").Value = "dateTime"
$ws.Range("H2").Value = "CodeableConcept"
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
